# Update the "Förändrad" date column (C) for rows 2-20 from 2023-09-13 (45182)
# to 2023-09-15 (45184), matching the commit's automatic data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 20; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value = 45184
    }
}
